$d = $word.ActiveDocument

$d.Content.Find.Execute("68×42=", $true, $false, $false, $false, $false, $true, 1, $false, "16×49=", 2) | Out-Null
$d.Content.Find.Execute("63×96=", $true, $false, $false, $false, $false, $true, 1, $false, "47×28=", 2) | Out-Null
$d.Content.Find.Execute("75×14=", $true, $false, $false, $false, $false, $true, 1, $false, "13×24=", 2) | Out-Null
$d.Content.Find.Execute("12×87=", $true, $false, $false, $false, $false, $true, 1, $false, "80×29=", 2) | Out-Null
$d.Content.Find.Execute("68×74=", $true, $false, $false, $false, $false, $true, 1, $false, "40×82=", 2) | Out-Null
$d.Content.Find.Execute("45×95=", $true, $false, $false, $false, $false, $true, 1, $false, "82×67=", 2) | Out-Null
$d.Content.Find.Execute("26×77=", $true, $false, $false, $false, $false, $true, 1, $false, "12×85=", 2) | Out-Null
$d.Content.Find.Execute("98×42=", $true, $false, $false, $false, $false, $true, 1, $false, "92×78=", 2) | Out-Null
$d.Content.Find.Execute("98×51=", $true, $false, $false, $false, $false, $true, 1, $false, "31×56=", 2) | Out-Null
$d.Content.Find.Execute("58×17=", $true, $false, $false, $false, $false, $true, 1, $false, "43×85=", 2) | Out-Null
$d.Content.Find.Execute("95×29=", $true, $false, $false, $false, $false, $true, 1, $false, "18×85=", 2) | Out-Null
$d.Content.Find.Execute("86×17=", $true, $false, $false, $false, $false, $true, 1, $false, "85×66=", 2) | Out-Null
$d.Content.Find.Execute("90×13=", $true, $false, $false, $false, $false, $true, 1, $false, "17×52=", 2) | Out-Null
$d.Content.Find.Execute("39×77=", $true, $false, $false, $false, $false, $true, 1, $false, "53×95=", 2) | Out-Null
$d.Content.Find.Execute("96×50=", $true, $false, $false, $false, $false, $true, 1, $false, "80×87=", 2) | Out-Null
$d.Content.Find.Execute("59×15=", $true, $false, $false, $false, $false, $true, 1, $false, "66×93=", 2) | Out-Null
$d.Content.Find.Execute("39×67=", $true, $false, $false, $false, $false, $true, 1, $false, "86×35=", 2) | Out-Null
$d.Content.Find.Execute("67×96=", $true, $false, $false, $false, $false, $true, 1, $false, "29×70=", 2) | Out-Null
$d.Content.Find.Execute("51×72=", $true, $false, $false, $false, $false, $true, 1, $false, "21×30=", 2) | Out-Null
$d.Content.Find.Execute("19×18=", $true, $false, $false, $false, $false, $true, 1, $false, "55×41=", 2) | Out-Null
$d.Content.Find.Execute("82×62=", $true, $false, $false, $false, $false, $true, 1, $false, "27×74=", 2) | Out-Null
$d.Content.Find.Execute("86×31=", $true, $false, $false, $false, $false, $true, 1, $false, "20×19=", 2) | Out-Null
$d.Content.Find.Execute("37×35=", $true, $false, $false, $false, $false, $true, 1, $false, "61×30=", 2) | Out-Null
$d.Content.Find.Execute("29×46=", $true, $false, $false, $false, $false, $true, 1, $false, "23×16=", 2) | Out-Null
$d.Content.Find.Execute("46×72=", $true, $false, $false, $false, $false, $true, 1, $false, "87×98=", 2) | Out-Null
